$d = $word.ActiveDocument

# Locate the paragraph that ends the "Things I Learned" list so far:
# "Very useful Collision Action Matrix of Documentation."
$anchorText = "Very useful Collision Action Matrix of Documentation."

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r`a") -eq $anchorText) {
        $target = $p
    }
}

if ($target -eq $null) {
    # Fallback: use the last paragraph in the body.
    $target = $d.Paragraphs.Last
}

# Insert a brand-new paragraph right after the anchor. Word automatically
# clones the paragraph formatting (pStyle "ListParagraph", numPr ilvl/numId,
# and rPr language) from the anchor paragraph onto the new one.
$target.Range.InsertParagraphAfter()

# Re-fetch the (now) last paragraph in the document -- this is the freshly
# inserted empty paragraph -- and give it its text.
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "Object Pool Design Pattern through Queue Data Structure."
